$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1804733727810651
$ws.Range("C2").Value = 0.5769230769230769
$ws.Range("J2").Value = 0.01479289940828402
$ws.Range("P2").Value = 0.1420118343195266
$ws.Range("S2").Value = 0.08579881656804733
$ws.Range("B3").Value = 0.005050505050505051
$ws.Range("C3").Value = 0.005050505050505051
$ws.Range("J3").Value = 0.01515151515151515
$ws.Range("P3").Value = 0.8080808080808081
$ws.Range("S3").Value = 0.1666666666666667
$ws.Range("J4").Value = 0.1186440677966102
$ws.Range("P4").Value = 0.559322033898305
$ws.Range("S4").Value = 0.3220338983050847
$ws.Range("B6").Value = 0.07804878048780488
$ws.Range("D6").Value = 0.01463414634146342
$ws.Range("F6").Value = 0.02926829268292683
$ws.Range("J6").Value = 0.2390243902439024
$ws.Range("O6").Value = 0.02926829268292683
$ws.Range("Q6").Value = 0.2097560975609756
$ws.Range("R6").Value = 0.09268292682926829
$ws.Range("S6").Value = 0.3073170731707317
$ws.Range("B7").Value = 0.12
$ws.Range("D7").Value = 0.01777777777777778
$ws.Range("E7").Value = 0.004444444444444444
$ws.Range("F7").Value = 0.04444444444444445
$ws.Range("J7").Value = 0.12
$ws.Range("O7").Value = 0.02222222222222222
$ws.Range("Q7").Value = 0.1733333333333333
$ws.Range("R7").Value = 0.04444444444444445
$ws.Range("S7").Value = 0.4533333333333333
$ws.Range("B8").Value = 0.1162280701754386
$ws.Range("D8").Value = 0.02192982456140351
$ws.Range("F8").Value = 0.04385964912280702
$ws.Range("J8").Value = 0.1184210526315789
$ws.Range("O8").Value = 0.02192982456140351
$ws.Range("Q8").Value = 0.2017543859649123
$ws.Range("R8").Value = 0.05921052631578947
$ws.Range("S8").Value = 0.4166666666666667
$ws.Range("B9").Value = 0.09333333333333334
$ws.Range("D9").Value = 0.02666666666666667
$ws.Range("F9").Value = 0.05333333333333334
$ws.Range("J9").Value = 0.1333333333333333
$ws.Range("O9").Value = 0.006666666666666667
$ws.Range("Q9").Value = 0.16
$ws.Range("R9").Value = 0.06666666666666667
$ws.Range("S9").Value = 0.46
$ws.Range("B10").Value = 0.1311343523732904
$ws.Range("D10").Value = 0.03137570394207562
$ws.Range("E10").Value = 0.0008045052292839903
$ws.Range("F10").Value = 0.06757843925985518
$ws.Range("J10").Value = 0.1086082059533387
$ws.Range("O10").Value = 0.01850362027353178
$ws.Range("Q10").Value = 0.2115848753016895
$ws.Range("R10").Value = 0.05309734513274336
$ws.Range("S10").Value = 0.3773129525341914
$ws.Range("G11").Value = 0.1585014409221902
$ws.Range("J11").Value = 0.1037463976945245
$ws.Range("K11").Value = 0.2161383285302594
$ws.Range("L11").Value = 0.5100864553314121
$ws.Range("S11").Value = 0.01152737752161383
$ws.Range("G12").Value = 0.7777777777777778
$ws.Range("J12").Value = 0.1722222222222222
$ws.Range("L12").Value = 0.02222222222222222
$ws.Range("S12").Value = 0.02777777777777778
$ws.Range("G13").Value = 0.6129032258064516
$ws.Range("J13").Value = 0.3870967741935484
$ws.Range("F15").Value = 0.02843601895734597
$ws.Range("H15").Value = 0.1658767772511848
$ws.Range("I15").Value = 0.04739336492890995
$ws.Range("J15").Value = 0.3744075829383886
$ws.Range("K15").Value = 0.06161137440758294
$ws.Range("M15").Value = 0.01421800947867299
$ws.Range("O15").Value = 0.03791469194312796
$ws.Range("S15").Value = 0.2701421800947867
$ws.Range("F16").Value = 0.04201680672268908
$ws.Range("H16").Value = 0.2016806722689076
$ws.Range("I16").Value = 0.04201680672268908
$ws.Range("J16").Value = 0.4033613445378151
$ws.Range("K16").Value = 0.1050420168067227
$ws.Range("M16").Value = 0.02100840336134454
$ws.Range("O16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.1134453781512605
$ws.Range("F17").Value = 0.01310043668122271
$ws.Range("H17").Value = 0.1746724890829694
$ws.Range("I17").Value = 0.08733624454148471
$ws.Range("J17").Value = 0.4235807860262009
$ws.Range("K17").Value = 0.1091703056768559
$ws.Range("M17").Value = 0.02838427947598253
$ws.Range("O17").Value = 0.05676855895196507
$ws.Range("S17").Value = 0.1069868995633188
$ws.Range("F18").Value = 0.007692307692307693
$ws.Range("H18").Value = 0.2153846153846154
$ws.Range("I18").Value = 0.06153846153846154
$ws.Range("J18").Value = 0.3692307692307693
$ws.Range("K18").Value = 0.1230769230769231
$ws.Range("M18").Value = 0.05384615384615385
$ws.Range("O18").Value = 0.09230769230769231
$ws.Range("S18").Value = 0.07692307692307693
$ws.Range("F19").Value = 0.02386634844868735
$ws.Range("H19").Value = 0.2084327764518695
$ws.Range("I19").Value = 0.06841686555290374
$ws.Range("J19").Value = 0.3627684964200477
$ws.Range("K19").Value = 0.1304693715194908
$ws.Range("M19").Value = 0.02784407319013524
$ws.Range("O19").Value = 0.066030230708035
$ws.Range("S19").Value = 0.1121718377088305
